$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the still-active "Great Sitkin" eruption (row 13): stop date and duration
# both advance by one day.
$ws.Range("C13").Value = 45492.0
$ws.Range("D13").Value = 1151.0

# Remove the rows that represented still-ongoing eruptions which have since
# been resolved/removed from the tracked data set. Deleting from the bottom
# up keeps the remaining row numbers stable while we work.
$ws.Range("A76:F76").EntireRow.Delete()
$ws.Range("A75:F75").EntireRow.Delete()
$ws.Range("A68:F68").EntireRow.Delete()
$ws.Range("A46:F46").EntireRow.Delete()
